# Fruta / hortaliza, semanal
# Insert a new record row at position 341 (shifting existing rows 341-389 down to 342-390)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 341; Excel shifts rows 341..389 down to 342..390
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new data record
$ws.Cells.Item(341, 1).Value = 3
$ws.Cells.Item(341, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(341, 3).Value = "Coquimbo"
$ws.Cells.Item(341, 4).Value = 44491
$ws.Cells.Item(341, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(341, 5).Value = 5
$ws.Cells.Item(341, 6).Value = 100112024
$ws.Cells.Item(341, 7).Value = "Choclo"
$ws.Cells.Item(341, 8).Value = "Dulce o Americano"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 100
$ws.Cells.Item(341, 11).Value = 41000
$ws.Cells.Item(341, 12).Value = 42000
$ws.Cells.Item(341, 13).Value = 41500
$ws.Cells.Item(341, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(341, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(341, 16).Value = 593
$ws.Cells.Item(341, 17).Value = 70
$ws.Cells.Item(341, 18).Value = "Hortaliza"
